# Updated 2D training schedules, no break screen
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New schedule data for rows 2-6 (A:J) — trialTrain, x_fixStart, y_fixStart,
# x_corrSteps, y_corrSteps, x_nrSteps, y_nrSteps, alienID, praclen, version
$data = @(
    @(1, 3, 7, 7, 5, 4, -2, 23, 5, "train_dim2_1"),
    @(2, 2, 7, 7, 6, 5, -1, 12, 5, "train_dim2_1"),
    @(3, 1, 8, 2, 3, 1, -5, 56, 5, "train_dim2_1"),
    @(4, 4, 9, 6, 5, 2, -4, 45, 5, "train_dim2_1"),
    @(5, 2, 5, 5, 2, 3, -3, 34, 5, "train_dim2_1")
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}

$ws.Range("I1").Select() | Out-Null
